$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to be treated as TEXT so Excel does not
    # auto-convert numeric-looking strings into numeric cells.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Drop the temporary "Text" number format again so the cell keeps
    # using the default style, while the stored value remains a string.
    $range.ClearFormats()
}

# Plain string equality/`-ceq` comparisons on worksheet names are
# case-insensitive in this environment, so "Vector_bf" and "Vector_BF"
# compare equal. Compare character codes explicitly to get a true
# case-sensitive comparison.
function Test-ExactNameMatch {
    param([string]$a, [string]$b)
    if ($a.Length -ne $b.Length) { return $false }
    $ca = $a.ToCharArray()
    $cb = $b.ToCharArray()
    for ($i = 0; $i -lt $ca.Length; $i++) {
        if ([int][char]$ca[$i] -ne [int][char]$cb[$i]) {
            return $false
        }
    }
    return $true
}

function Get-SheetByExactName {
    param([string]$name)
    foreach ($candidate in $wb.Worksheets) {
        if (Test-ExactNameMatch $candidate.Name $name) {
            return $candidate
        }
    }
    throw "Worksheet '$name' not found"
}

# --- Restricciones_del_follower -------------------------------------------------
$wsFollower = Get-SheetByExactName "Restricciones_del_follower"

Set-TextValue $wsFollower.Range("A2") "11.3 - 2x_1 + y_1 - y_2"
Set-TextValue $wsFollower.Range("B2") "-8.8"
Set-TextValue $wsFollower.Range("D2") "0.79"
Set-TextValue $wsFollower.Range("E2") "7.3"
Set-TextValue $wsFollower.Range("F2") "4.1"

Set-TextValue $wsFollower.Range("A3") "-3.55 + x_1 - 3x_2 + y_2"
Set-TextValue $wsFollower.Range("B3") "1.5499999999999998"
Set-TextValue $wsFollower.Range("D3") "0.09"
Set-TextValue $wsFollower.Range("E3") "2.2"
Set-TextValue $wsFollower.Range("F3") "3.7"

Set-TextValue $wsFollower.Range("A4") "-9.36 + x_1 + x_2"
Set-TextValue $wsFollower.Range("B4") "6.75"
Set-TextValue $wsFollower.Range("D4") "0.54"
Set-TextValue $wsFollower.Range("E4") "6.0"
Set-TextValue $wsFollower.Range("F4") "0.8999999999999999"

# --- Punto_modificado -------------------------------------------------
$wsPunto = Get-SheetByExactName "Punto_modificado"

Set-TextValue $wsPunto.Range("A2") "6.65"
Set-TextValue $wsPunto.Range("B2") "2.1"
Set-TextValue $wsPunto.Range("C2") "5.2"
Set-TextValue $wsPunto.Range("D2") "3.2"

# --- Vector_bf -------------------------------------------------
$wsBf = Get-SheetByExactName "Vector_bf"

Set-TextValue $wsBf.Range("A2") "3.21"
Set-TextValue $wsBf.Range("A3") "-0.29999999999999993"

# --- Vector_BF -------------------------------------------------
$wsBF = Get-SheetByExactName "Vector_BF"

Set-TextValue $wsBF.Range("A2") "8.399999999999999"
Set-TextValue $wsBF.Range("A3") "-0.39999999999999947"
Set-TextValue $wsBF.Range("A4") "-7.8"
Set-TextValue $wsBF.Range("A5") "5.1"
